# Re-pull / push updated "dSF" (column F) data for each row of game log data.
# The sheet layout (row -> column) is:
#   B:date  C:TB  D:PC  E:dS0  F:dSF  G:K  H:IP  I:I0  J:IF
# Only column F ("dSF") changes in this pass; row 9 already matches the
# freshly pulled value and is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dSF = @{
    2  = -1
    3  = -3
    4  = 6
    5  = -2
    6  = 1
    7  = -1
    8  = 1
    10 = 4
    11 = -2
    12 = 3
    13 = 2
    14 = -2
    15 = 1
    16 = -3
    17 = 1
    18 = 3
    19 = -4
    20 = -1
    21 = -3
    22 = 5
    23 = -2
    24 = 3
    25 = 9
    26 = 6
    27 = 1
    28 = -1
    29 = -1
}

foreach ($row in $dSF.Keys) {
    $ws.Cells.Item($row, 6).Value = $dSF[$row]
}
